# Figure 1 label updates: rename the "Flood" treatment labels to
# "Waterlogged" labels in the treatment column (column C) of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Range("C1:C80")
$col.Replace("Flood24", "Waterlogged 24")
$col.Replace("Flood72", "Waterlogged 72")
$col.Replace("Flood48", "Waterlogged 48")

# Restore the last-used selection recorded in the saved workbook.
$ws.Range("G17").Select()
